$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (interested count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1191
$ws1.Range("F4").Value = 2662
$ws1.Range("F5").Value = 234

# Sheet "全部类型" - update "想去人数" (interested count) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1191
$ws4.Range("F6").Value = 2662
$ws4.Range("F8").Value = 234
